$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 141, shifting existing rows 141:271 down to 142:272
$ws.Rows("141:141").Insert()

# Populate the newly inserted row 141 with the new weekly price record
$ws.Range("A141").Value = 3
$ws.Range("B141").Value = "Femacal de La Calera"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44587
$ws.Range("E141").Value = 5
$ws.Range("F141").Value = 100112012
$ws.Range("G141").Value = "Espinaca"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 60
$ws.Range("K141").Value = 4000
$ws.Range("L141").Value = 4000
$ws.Range("M141").Value = 4000
$ws.Range("N141").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O141").Value = "Provincia de Quillota"
$ws.Range("P141").Value = 1333
$ws.Range("Q141").Value = 3
$ws.Range("R141").Value = "Hortaliza"
